# Swap the contents of rows 17 and 18 (record data for two observations
# that were reordered), for the columns that actually differ between the
# two records: A, B, E, F, G, H, M, Q, R.
# (Columns D, I, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY are
# identical between the two rows, so they are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the current ("before") values of row 17 and row 18 ---
# Note: use Value2 (not Value) to get plain scalar data back from the
# COM property getter in this runtime.
$A17 = $ws.Range("A17").Value2
$B17 = $ws.Range("B17").Value2
$E17 = $ws.Range("E17").Value2
$F17 = $ws.Range("F17").Value2
$G17 = $ws.Range("G17").Value2
$H17 = $ws.Range("H17").Value2
$M17 = $ws.Range("M17").Value2
$Q17 = $ws.Range("Q17").Value2
$R17 = $ws.Range("R17").Value2

$A18 = $ws.Range("A18").Value2
$B18 = $ws.Range("B18").Value2
$E18 = $ws.Range("E18").Value2
$F18 = $ws.Range("F18").Value2
$G18 = $ws.Range("G18").Value2
$H18 = $ws.Range("H18").Value2
$M18 = $ws.Range("M18").Value2
$Q18 = $ws.Range("Q18").Value2
$R18 = $ws.Range("R18").Value2

# --- Write row 17 with the values previously in row 18 ---
$ws.Range("A17").Value2 = $A18
$ws.Range("B17").Value2 = $B18
$ws.Range("E17").Value2 = $E18
$ws.Range("F17").Value2 = $F18
$ws.Range("G17").Value2 = $G18
$ws.Range("H17").Value2 = $H18
$ws.Range("M17").Value2 = $M18
$ws.Range("Q17").Value2 = $Q18
$ws.Range("R17").Value2 = $R18

# --- Write row 18 with the values previously in row 17 ---
$ws.Range("A18").Value2 = $A17
$ws.Range("B18").Value2 = $B17
$ws.Range("E18").Value2 = $E17
$ws.Range("F18").Value2 = $F17
$ws.Range("G18").Value2 = $G17
$ws.Range("H18").Value2 = $H17
$ws.Range("M18").Value2 = $M17
$ws.Range("Q18").Value2 = $Q17
$ws.Range("R18").Value2 = $R17
